$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.Style = 'Normal'
}

$ws.Range('D2').Value = '63.321.65'
$ws.Range('E2').Value = '  +5.09%  '

$ws.Range('D3').Value = '2.712.96'
$ws.Range('E3').Value = '  +4.29%  '

$ws.Range('E4').Value = '  +0.00%  '

Set-TextValue 'D5' '586.74'
$ws.Range('E5').Value = '  +0.36%  '

Set-TextValue 'D6' '149.87'
$ws.Range('E6').Value = '  +4.91%  '

$ws.Range('E7').Value = '  -0.35%  '

Set-TextValue 'D8' '0.607'
$ws.Range('E8').Value = '  +1.22%  '

$ws.Range('D9').Value = '2.742.54'
$ws.Range('E9').Value = '  +5.20%  '

Set-TextValue 'D10' '6.72'
$ws.Range('E10').Value = '  +2.74%  '

$ws.Range('E11').Value = '  +7.61%  '

$ws.Range('E12').Value = '  +5.04%  '

$ws.Range('D14').Value = '3.188.16'
$ws.Range('E14').Value = '  +4.15%  '

Set-TextValue 'D15' '26.60'
$ws.Range('E15').Value = '  +8.63%  '

$ws.Range('D16').Value = '63.203.69'
$ws.Range('E16').Value = '  +4.89%  '

$ws.Range('E17').Value = '  +7.46%  '

$ws.Range('D18').Value = '2.737.23'
$ws.Range('E18').Value = '  +5.13%  '

Set-TextValue 'D19' '11.93'
$ws.Range('E19').Value = '  +5.30%  '

Set-TextValue 'D20' '4.86'
$ws.Range('E20').Value = '  +5.29%  '

Set-TextValue 'D21' '364.27'
$ws.Range('E21').Value = '  +5.28%  '

Set-TextValue 'D22' '6.98'
$ws.Range('E22').Value = '  +1.12%  '

$ws.Range('E23').Value = '  -0.36%  '

$ws.Range('E24').Value = '  +0.59%  '

Set-TextValue 'D25' '65.54'
$ws.Range('E25').Value = '  +2.77%  '

$ws.Range('E26').Value = '  +4.04%  '

Set-TextValue 'D27' '8.54'
$ws.Range('E27').Value = '  +6.90%  '

$ws.Range('E28').Value = '  -0.23%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0860'
$ws.Range('E29').Value = '  +7.92%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '2.03'
$ws.Range('E30').Value = '  +6.56%  '

Set-TextValue 'D31' '7.07'
$ws.Range('E31').Value = '  +10.52%  '

Set-TextValue 'D32' '170.89'
$ws.Range('E32').Value = '  +1.52%  '

$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D33' '0.996'
$ws.Range('E33').Value = '  -0.19%  '

$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D34' '1.19'
$ws.Range('E34').Value = '  +21.22%  '

Set-TextValue 'D35' '20.50'
$ws.Range('E35').Value = '  +5.57%  '

Set-TextValue 'D36' '4.76'
$ws.Range('E36').Value = '  +11.98%  '

$ws.Range('E37').Value = '  +8.73%  '

$ws.Range('E38').Value = '  +10.68%  '

$ws.Range('E39').Value = '  +19.38%  '

Set-TextValue 'D40' '352.44'
$ws.Range('E40').Value = '  +12.13%  '

$ws.Range('E41').Value = '  +9.80%  '

Set-TextValue 'D42' '39.20'
$ws.Range('E42').Value = '  +2.67%  '

Set-TextValue 'D43' '5.65'
$ws.Range('E43').Value = '  +13.53%  '

Set-TextValue 'D44' '21.63'
$ws.Range('E44').Value = '  +9.06%  '

$ws.Range('E45').Value = '  +7.75%  '

Set-TextValue 'D46' '21.68'
$ws.Range('E46').Value = '  +9.15%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0260'
$ws.Range('E47').Value = '  +7.16%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D48' '138.40'
$ws.Range('E48').Value = '  +2.01%  '

Set-TextValue 'D49' '0.641'
$ws.Range('E49').Value = '  +5.85%  '

$ws.Range('E50').Value = '  +1.80%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.145.13'
$ws.Range('E51').Value = '  +6.38%  '
